$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.724.33"
$ws.Range("E2").Value = "  -1.64%  "
$ws.Range("D3").Value = "2.528.57"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.15"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.77"
$ws.Range("E6").Value = "  +3.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.31"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0804"
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "2.927.35"
$ws.Range("E14").Value = "  -2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.66"
$ws.Range("E15").Value = "  +2.36%  "
$ws.Range("D16").Value = "2.484.82"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("E17").Value = "  -4.66%  "
$ws.Range("D18").Value = "42.694.84"
$ws.Range("E18").Value = "  -1.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  -2.41%  "
$ws.Range("D20").Value = "0.0₃0951"
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.24"
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.44"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.69"
$ws.Range("E23").Value = "  -4.20%  "
$ws.Range("E24").Value = "  -3.31%  "
$ws.Range("E25").Value = "  -2.17%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.08"
$ws.Range("E27").Value = "  -5.40%  "
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.15"
$ws.Range("E29").Value = "  -5.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.18"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.33"
$ws.Range("E32").Value = "  -0.82%  "
$ws.Range("E33").Value = "  +8.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0793"
$ws.Range("E34").Value = "  -2.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.63"
$ws.Range("E35").Value = "  -2.74%  "
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.19"
$ws.Range("E37").Value = "  -9.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.20"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.33"
$ws.Range("E41").Value = "  +6.16%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0300"
$ws.Range("E45").Value = "  -1.79%  "
$ws.Range("D46").Value = "1.986.39"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.88"
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("D48").Value = "2.782.27"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "80.68"
$ws.Range("E49").Value = "  -4.08%  "
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.850"
$ws.Range("E51").Value = "  +7.44%  "
